# Refresh cached Market Board figures (currentAveragePrice*, Leve cost/profit
# columns H:N) on a handful of leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/WVR
# sheets, per the latest scrape. One block per edited row.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 496.57144
$ws.Range("I18").Value = 246
$ws.Range("K18").Value = 246
$ws.Range("M18").Value = 38

$ws.Range("H98").Value = 4465.0347
$ws.Range("I98").Value = 4465.0347
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 4465.0347
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -2967.0347
$ws.Range("N98").ClearContents()

$ws.Range("H122").Value = 4465.0347
$ws.Range("I122").Value = 4465.0347
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 13395.1041
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -10945.1041
$ws.Range("N122").ClearContents()

$ws.Range("H137").Value = 2266.475
$ws.Range("I137").Value = 2018.3889
$ws.Range("K137").Value = 6055.1667
$ws.Range("M137").Value = -3505.1667

$ws.Range("H138").Value = 2430.859
$ws.Range("J138").Value = 2520.4768
$ws.Range("L138").Value = 7561.430399999999
$ws.Range("N138").Value = -17841.4304

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()

$ws.Range("H32").Value = 9881.951999999999
$ws.Range("I32").Value = 7308.507
$ws.Range("K32").Value = 7308.507
$ws.Range("M32").Value = -7021.507

$ws.Range("H45").Value = 901.8889
$ws.Range("I45").Value = 788.93335
$ws.Range("K45").Value = 788.93335
$ws.Range("M45").Value = -411.93335

$ws.Range("H132").Value = 3785.2334
$ws.Range("I132").Value = 3185.85
$ws.Range("K132").Value = 9557.549999999999
$ws.Range("M132").Value = -7027.549999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3151.111
$ws.Range("J20").Value = 2820
$ws.Range("L20").Value = 2820
$ws.Range("N20").Value = -3314

$ws.Range("H81").Value = 19751.4
$ws.Range("J81").Value = 19751.4
$ws.Range("L81").Value = 19751.4
$ws.Range("N81").Value = -21873.4

$ws.Range("H84").Value = 19751.4
$ws.Range("J84").Value = 19751.4
$ws.Range("L84").Value = 59254.2
$ws.Range("N84").Value = -69862.20000000001

$ws.Range("H105").Value = 111123380
$ws.Range("I105").Value = 125013650
$ws.Range("K105").Value = 125013650
$ws.Range("M105").Value = -125011903

$ws.Range("H107").Value = 1018.5238
$ws.Range("I107").Value = 886.8333
$ws.Range("K107").Value = 886.8333
$ws.Range("M107").Value = 1033.1667

$ws.Range("H134").Value = 18176.5
$ws.Range("I134").Value = 1765
$ws.Range("K134").Value = 5295
$ws.Range("M134").Value = -2760

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1622.0952
$ws.Range("I31").Value = 1588.4878
$ws.Range("K31").Value = 1588.4878
$ws.Range("M31").Value = -1293.4878

$ws.Range("H34").Value = 1622.0952
$ws.Range("I34").Value = 1588.4878
$ws.Range("K34").Value = 1588.4878
$ws.Range("M34").Value = -1386.4878

$ws.Range("H64").Value = 20000
$ws.Range("J64").Value = 20000
$ws.Range("L64").Value = 20000
$ws.Range("N64").Value = -20496

$ws.Range("H67").Value = 20000
$ws.Range("J67").Value = 20000
$ws.Range("L67").Value = 20000
$ws.Range("N67").Value = -21716

$ws.Range("H99").Value = 1802.0714
$ws.Range("I99").Value = 1682.9
$ws.Range("J99").Value = 2100
$ws.Range("K99").Value = 1682.9
$ws.Range("L99").Value = 2100
$ws.Range("M99").Value = -184.9000000000001
$ws.Range("N99").Value = -5096

$ws.Range("H105").Value = 836.125
$ws.Range("I105").Value = 812.7143
$ws.Range("K105").Value = 812.7143
$ws.Range("M105").Value = 934.2857

$ws.Range("H107").Value = 1156.7778
$ws.Range("I107").Value = 451.83334
$ws.Range("J107").Value = 2566.6667
$ws.Range("K107").Value = 451.83334
$ws.Range("L107").Value = 2566.6667
$ws.Range("M107").Value = 1468.16666
$ws.Range("N107").Value = -6406.6667

$ws.Range("H126").Value = 1802.0714
$ws.Range("I126").Value = 1682.9
$ws.Range("J126").Value = 2100
$ws.Range("K126").Value = 5048.700000000001
$ws.Range("L126").Value = 6300
$ws.Range("M126").Value = -2578.700000000001
$ws.Range("N126").Value = -11240

$ws.Range("H141").Value = 1213598
$ws.Range("J141").Value = 1213598
$ws.Range("L141").Value = 1213598
$ws.Range("N141").Value = -1223958

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 3617.5881
$ws.Range("I81").Value = 2599.75
$ws.Range("J81").Value = 3930.7693
$ws.Range("K81").Value = 7799.25
$ws.Range("L81").Value = 11792.3079
$ws.Range("M81").Value = -6676.25
$ws.Range("N81").Value = -14038.3079

$ws.Range("H84").Value = 3617.5881
$ws.Range("I84").Value = 2599.75
$ws.Range("J84").Value = 3930.7693
$ws.Range("K84").Value = 23397.75
$ws.Range("L84").Value = 35376.9237
$ws.Range("M84").Value = -17781.75
$ws.Range("N84").Value = -46608.9237

$ws.Range("H131").Value = 21772650
$ws.Range("I131").Value = 71429140
$ws.Range("J131").Value = 47937.344
$ws.Range("K131").Value = 214287420
$ws.Range("L131").Value = 143812.032
$ws.Range("M131").Value = -214282380
$ws.Range("N131").Value = -153892.032

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 20000
$ws.Range("J45").Value = 20000
$ws.Range("L45").Value = 20000
$ws.Range("N45").Value = -21118

$ws.Range("H132").Value = 7165.25
$ws.Range("I132").Value = 14505.5
$ws.Range("J132").Value = 4718.5
$ws.Range("K132").Value = 43516.5
$ws.Range("L132").Value = 14155.5
$ws.Range("M132").Value = -40986.5
$ws.Range("N132").Value = -19215.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 15900
$ws.Range("J74").Value = 15900
$ws.Range("L74").Value = 15900
$ws.Range("N74").Value = -17772

$ws.Range("H77").Value = 15900
$ws.Range("J77").Value = 15900
$ws.Range("L77").Value = 47700
$ws.Range("N77").Value = -57060

$ws.Range("H100").Value = 325.75
$ws.Range("I100").Value = 325.75
$ws.Range("K100").Value = 651.5
$ws.Range("M100").Value = -110.5

$ws.Range("H107").Value = 542.25
$ws.Range("I107").Value = 489.83334
$ws.Range("J107").Value = 699.5
$ws.Range("K107").Value = 1469.50002
$ws.Range("L107").Value = 2098.5
$ws.Range("M107").Value = 450.4999800000001
$ws.Range("N107").Value = -5938.5

$ws.Range("H136").Value = 1669.4
$ws.Range("I136").Value = 1521.5555
$ws.Range("K136").Value = 4564.666499999999
$ws.Range("M136").Value = -2014.666499999999
